$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 25; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains(" / ")) {
        $cell.Value2 = $val.Replace(" / ", "/")
    }
}
